# Updates the "Price" (column D) and "Volume(1h)" (column E) cells of the
# crypto tracker sheet with refreshed quote data, mirroring the scheduled
# GitHub Actions data-refresh commit ("Updated cryptos list ...").
#
# Each touched cell keeps its original plain-text storage: the price/volume
# strings (e.g. "61.281.93", "  +0.17%  ") are not valid numeric literals
# for some rows but ARE for others (e.g. "596.54"), and Excel's normal
# Range.Value assignment auto-converts number-looking text into a real
# number. Forcing the NumberFormat to Text ("@") before the write, then
# resetting the style back to "Normal" afterwards, keeps every value a
# plain string (matching the source data) without leaving a stray
# per-cell number format behind.

function Set-CellText($sheet, $ref, $val) {
    $cell = $sheet.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "61.281.93"
Set-CellText $ws "E2" "  +0.17%  "
Set-CellText $ws "D3" "2.923.42"
Set-CellText $ws "E3" "  -0.29%  "
Set-CellText $ws "E4" "  +0.01%  "
Set-CellText $ws "D5" "596.54"
Set-CellText $ws "E5" "  +0.43%  "
Set-CellText $ws "D6" "144.98"
Set-CellText $ws "E6" "  -1.00%  "
Set-CellText $ws "E7" "  +0.03%  "
Set-CellText $ws "D8" "0.501"
Set-CellText $ws "E8" "  -1.15%  "
Set-CellText $ws "D9" "6.96"
Set-CellText $ws "E9" "  +0.70%  "
Set-CellText $ws "E10" "  -2.73%  "
Set-CellText $ws "E11" "  -0.87%  "
Set-CellText $ws "E12" "  -1.51%  "
Set-CellText $ws "D13" "33.41"
Set-CellText $ws "E13" "  -1.32%  "
Set-CellText $ws "E14" "  +0.27%  "
Set-CellText $ws "D15" "3.410.60"
Set-CellText $ws "E15" "  -0.14%  "
Set-CellText $ws "D16" "61.400.79"
Set-CellText $ws "E16" "  +0.42%  "
Set-CellText $ws "D17" "2.922.96"
Set-CellText $ws "E17" "  -0.33%  "
Set-CellText $ws "D18" "6.68"
Set-CellText $ws "E18" "  -0.75%  "
Set-CellText $ws "D19" "430.87"
Set-CellText $ws "E19" "  -0.32%  "
Set-CellText $ws "D20" "13.45"
Set-CellText $ws "E20" "  -0.20%  "
Set-CellText $ws "D21" "0.675"
Set-CellText $ws "E21" "  -1.37%  "
Set-CellText $ws "D22" "7.05"
Set-CellText $ws "E22" "  -0.73%  "
Set-CellText $ws "D23" "81.78"
Set-CellText $ws "E23" "  +0.43%  "
Set-CellText $ws "D24" "10.83"
Set-CellText $ws "E24" "  -2.19%  "
Set-CellText $ws "E25" "  -2.57%  "
Set-CellText $ws "D26" "11.72"
Set-CellText $ws "E26" "  -2.73%  "
Set-CellText $ws "E27" "  +0.04%  "
Set-CellText $ws "E28" "  -5.40%  "
Set-CellText $ws "E29" "  -0.56%  "
Set-CellText $ws "D30" "6.89"
Set-CellText $ws "E30" "  -3.25%  "
Set-CellText $ws "D31" "0.109"
Set-CellText $ws "E31" "  +1.14%  "
Set-CellText $ws "D32" "26.55"
Set-CellText $ws "E32" "  +0.14%  "
Set-CellText $ws "E33" "  +0.08%  "
Set-CellText $ws "D34" "0.0₃0881"
Set-CellText $ws "E34" "  +2.32%  "
Set-CellText $ws "E35" "  -0.20%  "
Set-CellText $ws "D36" "5.61"
Set-CellText $ws "E36" "  -0.68%  "
Set-CellText $ws "D37" "2.97"
Set-CellText $ws "E37" "  -3.73%  "
Set-CellText $ws "D38" "1.99"
Set-CellText $ws "E38" "  -0.57%  "
Set-CellText $ws "D39" "0.122"
Set-CellText $ws "E39" "  -2.11%  "
Set-CellText $ws "D40" "8.55"
Set-CellText $ws "E40" "  -0.78%  "
Set-CellText $ws "D41" "42.17"
Set-CellText $ws "E41" "  +6.28%  "
Set-CellText $ws "D42" "0.281"
Set-CellText $ws "E42" "  -2.64%  "
Set-CellText $ws "D43" "0.0345"
Set-CellText $ws "E43" "  -0.53%  "
Set-CellText $ws "D44" "2.697.66"
Set-CellText $ws "E44" "  -0.89%  "
Set-CellText $ws "D45" "133.71"
Set-CellText $ws "E45" "  +2.30%  "
Set-CellText $ws "D46" "361.85"
Set-CellText $ws "E46" "  -4.11%  "
Set-CellText $ws "E47" "  +0.03%  "
Set-CellText $ws "D48" "23.47"
Set-CellText $ws "E48" "  -3.06%  "
Set-CellText $ws "E49" "  -1.53%  "
Set-CellText $ws "D50" "1.99"
Set-CellText $ws "E50" "  -2.23%  "
Set-CellText $ws "E51" "  -2.71%  "
